# Append a new row 44 (copy of row 43, with the hour bumped by 1) to each
# of the 4 worksheets in the workbook, matching the sensor-log pattern
# already present in each sheet.

$wb = $excel.ActiveWorkbook

# Per-sheet data for the new row 44: A(time), B, C, D, E, F, G, H, I
$rowsData = @{
    "ROW35-FE-LIFTER"  = @("2025-03-06 03:42:06", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,", "0x01,0x90,", "0x d", 400, "568631262647113770877196", 400, 13)
    "ROW35-MID-LIFTER" = @("2025-03-06 03:29:35", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,", "0x01,0x90,", "0x e", 400, "568631262647113770942732", 400, 14)
    "ROW02-FE-LIFTER"  = @("2025-03-06 03:51:45", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,", "0x01,0x90,", "0xff", 400, "568631262647113769959692", 400, 255)
    "ROW02-MID-LIFTER" = @("2025-03-06 03:41:15", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x01,0x90,", "0x 3", 400, "568631262647113769959692", 400, 3)
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($rowsData.ContainsKey($name)) {
        $vals = $rowsData[$name]
        $r = 44

        $ws.Cells.Item($r, 1).Value = $vals[0]   # A: time (text)
        $ws.Cells.Item($r, 2).Value = $vals[1]   # B: total length hex (text)
        $ws.Cells.Item($r, 3).Value = $vals[2]   # C: ID hex (text)
        $ws.Cells.Item($r, 4).Value = $vals[3]   # D: actual length hex (text)
        $ws.Cells.Item($r, 5).Value = $vals[4]   # E: checksum hex (text)
        $ws.Cells.Item($r, 6).Value = $vals[5]   # F: total length dec (number)

        # G holds a 24-digit integer that exceeds floating point precision,
        # so it must be stored as text, not auto-coerced to a number.
        # Force text format, assign, then clear the format again so no
        # leftover style index is left on the cell (matches the rest of
        # the column, which carries no cell style).
        $ws.Cells.Item($r, 7).NumberFormat = "@"
        $ws.Cells.Item($r, 7).Value = $vals[6]   # G: ID dec (text, big integer)
        $ws.Cells.Item($r, 7).ClearFormats()

        $ws.Cells.Item($r, 8).Value = $vals[7]   # H: actual length dec (number)
        $ws.Cells.Item($r, 9).Value = $vals[8]   # I: checksum dec (number)
    }
}
